$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared string values reused from existing rows in the sheet
$subjectIds   = "1-7 (training), 0 (testing)"
$textEncoder  = "Default CLIP"
$imageEncoder = "3D Resnet18"
$similarity   = "Cosine Similarity"
$clipHyper    = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers"
$optHyper     = "LR=1e-5, batch_size=32, weight_decay=0.2"
$trainValTest = "700/0/100"
$imageInput   = "(fmri channel for each word) (detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word"
$textTokenizer = "Default CLIP"
$textInput4   = "4 words"
$textInput8   = "8 words"

# Row 55: fill in the rest of the row that previously only had C55 set
$ws.Range("A55").Value = $subjectIds
$ws.Range("B55").Value = $textInput4
$ws.Range("C55").Value = $imageInput
$ws.Range("D55").Value = $textTokenizer
$ws.Range("E55").Value = $textEncoder
$ws.Range("F55").Value = $imageEncoder
$ws.Range("G55").Value = $similarity
$ws.Range("H55").Value = $clipHyper
$ws.Range("I55").Value = $optHyper
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = $trainValTest

# Row 56: new row with same data, but "8 words" text input
$ws.Range("A56").Value = $subjectIds
$ws.Range("B56").Value = $textInput8
$ws.Range("C56").Value = $imageInput
$ws.Range("D56").Value = $textTokenizer
$ws.Range("E56").Value = $textEncoder
$ws.Range("F56").Value = $imageEncoder
$ws.Range("G56").Value = $similarity
$ws.Range("H56").Value = $clipHyper
$ws.Range("I56").Value = $optHyper
$ws.Range("J56").Value = 50
$ws.Range("K56").Value = $trainValTest

# Reflect updated selection/view similar to authored workbook
$ws.Range("P56").Select()
